# Commit: Wed, May 06, 2020  4:10:22 AM
#
# 1) The table on slide 5 gets its table style (tableStyleId) changed
#    from {39B57D46-C57B-4F33-9AA4-25520CA80B51} to
#    {5362F531-261E-42E8-9F1D-53CBD097AA9B}.
# 2) The presentation's applied colour theme is switched from the
#    "Integral" / "Red Violet" palette back to the standard Office
#    theme palette (dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink).

$p = $ppt.ActivePresentation

# --- 1. Re-style the table on the slide that has one ------------------
$oldStyleId = "{39B57D46-C57B-4F33-9AA4-25520CA80B51}"
$newStyleId = "{5362F531-261E-42E8-9F1D-53CBD097AA9B}"

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shp = $slide.Shapes.Item($shi)
        if ($shp.HasTable) {
            $tbl = $shp.Table
            if ($tbl.Style -eq $oldStyleId) {
                $tbl.ApplyStyle($newStyleId)
            }
        }
    }
}

# --- 2. Swap the active theme's colour scheme back to "Office" --------
# dk1, lt1, dk2, lt2, accent1, accent2, accent3, accent4, accent5,
# accent6, hlink, folHlink (MsoThemeColorSchemeIndex order 1-12).
$officeThemeHex = @(
    "000000", "FFFFFF", "44546A", "E7E6E6",
    "5B9BD5", "ED7D31", "A5A5A5", "FFC000",
    "4472C4", "70AD47", "0563C1", "954F72"
)

$slideOne = $p.Slides.Item(1)
$colorScheme = $slideOne.ThemeColorScheme

for ($i = 1; $i -le $colorScheme.Count; $i++) {
    $hex = $officeThemeHex[$i - 1]
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    $colorScheme.Item($i).RGB = $r + ($g * 256) + ($b * 65536)
}
